$wb = $excel.ActiveWorkbook

# --- Schedule sheet ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 772.50727125
$schedule.Range("F2").Value = 12.77293768601191

# --- Detailed sheet ---
$detailed = $wb.Worksheets.Item("Detailed")
$detailed.Range("B5").Value = 36.2
$detailed.Range("B6").Value = 23.91419
$detailed.Range("B7").Value = 57.06003
$detailed.Range("B8").Value = 57.06003
$detailed.Range("B9").Value = 57.83179
$detailed.Range("B10").Value = 57.48786
$detailed.Range("B11").Value = 60.39169
$detailed.Range("B12").Value = 60.34276
$detailed.Range("B13").Value = 71.95779
$detailed.Range("B14").Value = 65
$detailed.Range("B16").Value = 35.88
$detailed.Range("B17").Value = 7.82154
$detailed.Range("B18").Value = 0.7
$detailed.Range("B19").Value = 7.00436
$detailed.Range("B20").Value = 0.7
$detailed.Range("B21").Value = -0.94872
$detailed.Range("B22").Value = -5.01
$detailed.Range("B23").Value = -5.01
$detailed.Range("B24").Value = 0
$detailed.Range("B25").Value = -5.82017
$detailed.Range("B26").Value = -4.658
$detailed.Range("B27").Value = -5.17224
$detailed.Range("B28").Value = -5.50985
$detailed.Range("B29").Value = -2.83936
$detailed.Range("B30").Value = -0.93813
$detailed.Range("B31").Value = 0.00002
$detailed.Range("B32").Value = 0.009469999999999999
$detailed.Range("B33").Value = 0.7
$detailed.Range("B34").Value = 1.21401
$detailed.Range("B35").Value = -0.91614
$detailed.Range("B36").Value = 0
$detailed.Range("B37").Value = 0.33937
$detailed.Range("B38").Value = 12.09882
$detailed.Range("B39").Value = 43.32325
$detailed.Range("B40").Value = 56.77711
$detailed.Range("B41").Value = 60.96077
$detailed.Range("B42").Value = 65
$detailed.Range("B44").Value = 62.69405
$detailed.Range("B45").Value = 62.96752
$detailed.Range("B46").Value = 61.4952
$detailed.Range("B47").Value = 58.97571
$detailed.Range("B48").Value = 58.48808
$detailed.Range("B49").Value = 62.01329

# Type column updates (forecast -> historical)
$detailed.Range("C7").Value = "historical"
$detailed.Range("C8").Value = "historical"
